# "Generate Report for Handback"
# The CI run produced a (too-old) handback for d91bbd79-3ecf-44ad-9225-acf5028c185e
# on both the zh-cn and de-de status sheets (row 7 of each sheet's table).
# Fill in: Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K) and Error Detail (P), and add the
# corresponding hyperlink on the new "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$mdName   = "d91bbd79-3ecf-44ad-9225-acf5028c185e.md"
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/291da1baa54b13ff2b9833d0449ee8be6f10359f/e2e/d91bbd79-3ecf-44ad-9225-acf5028c185e.md"
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cab8c1ed2c7609976eee60039ac05ba23999443a/e2e/d91bbd79-3ecf-44ad-9225-acf5028c185e.md"
$errorDetail = "The version of handback file is not the latest, current: " + $currentUrl + ", latest: " + $latestUrl + "."

# zh-cn sheet
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("I7").Value = $mdName
$ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", $mdName) | Out-Null
$ws.Range("J7").Value = "d91bbd79-3ecf-44ad-9225-acf5028c185e.79f2d38d59f839531f6c26f4b1a927a4ccc7b276.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-23 00:52:46"
$ws.Range("P7").Value = $errorDetail

# de-de sheet
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("I7").Value = $mdName
$ws2.Hyperlinks.Add($ws2.Range("I7"), $latestUrl, "", "", $mdName) | Out-Null
$ws2.Range("J7").Value = "d91bbd79-3ecf-44ad-9225-acf5028c185e.79f2d38d59f839531f6c26f4b1a927a4ccc7b276.de-de.xlf"
$ws2.Range("K7").Value = "2016-08-23 00:52:53"
$ws2.Range("P7").Value = $errorDetail
